$d = $word.ActiveDocument

# -------------------------------------------------------------------------
# 1) Formula "(TI - TO) / DT" -> "(TO - TI) / DT" (swap subscripted I/O and
#    change the separating hyphen to an en dash), in the "SHT = HT (OD/TC)..."
#    paragraph. We locate the paragraph via its distinctive leading text and
#    edit the individual characters so the subscript formatting on the I/O
#    letters is preserved.
# -------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $ptext = $para.Range.Text
    if ($ptext.StartsWith("`tSHT`t= HT")) {
        $p = $para.Range
        $pstart = $p.Start

        # Find the index of the "(TI" sequence within the paragraph text so this
        # is resilient to any minor offset differences.
        $idx = $ptext.IndexOf("(TI")
        # $idx is 0-based; characters are 1-based, so the "T" of "TI" is at
        # position $idx + 2 (1-based) and "I" immediately follows it.
        $posT1 = $idx + 2   # 1-based index of first "T"
        $posI  = $idx + 3   # 1-based index of "I"
        $posSp1 = $idx + 4  # space after "I"
        $posHy  = $idx + 5  # hyphen
        $posSp2 = $idx + 6  # subscript space
        $posT2  = $idx + 7  # second "T"
        $posO   = $idx + 8  # "O"

        # Swap the subscripted letters, preserving their run formatting.
        $p.Characters($posI).Text = "O"
        $p.Characters($posO).Text = "I"

        # Replace " -" (space + hyphen) with " " + en dash, preserving the
        # run/formatting of the characters around it.
        $absStart = $pstart + ($posSp1 - 1)
        $absEnd = $pstart + $posHy
        $dashRange = $d.Range($absStart, $absEnd)
        $dashRange.Text = " " + [char]0x2013
        break
    }
}

# -------------------------------------------------------------------------
# 2) Grammar-check punctuation fixes: "label; value" / "label, value" -> 
#    "label: value" (semicolons/commas used before an explanatory value
#    become colons).
# -------------------------------------------------------------------------
$d.Content.Find.Execute("Conversion constant; 6", $true, $false, $false, $false, $false, $true, 1, $false, "Conversion constant: 6", 2) | Out-Null
$d.Content.Find.Execute("peak demand; 100%", $true, $false, $false, $false, $false, $true, 1, $false, "peak demand: 100%", 2) | Out-Null
$d.Content.Find.Execute("HVAC system; `${OHS}", $true, $false, $false, $false, $false, $true, 1, $false, "HVAC system: `${OHS}", 2) | Out-Null
$d.Content.Find.Execute("Conversion constant, 0.746", $true, $false, $false, $false, $false, $true, 1, $false, "Conversion constant: 0.746", 2) | Out-Null
$d.Content.Find.Execute("air curtain motors; `${HP}", $true, $false, $false, $false, $false, $true, 1, $false, "air curtain motors: `${HP}", 2) | Out-Null
$d.Content.Find.Execute("air curtains, `${OHAC}", $true, $false, $false, $false, $false, $true, 1, $false, "air curtains: `${OHAC}", 2) | Out-Null
